$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Common data describing the 5 "File Name" entities used across the 3 sheets.
# Each entry: guid, content-hash (used in xlf file names), zh-cn handoff
# datetime, de-de handoff datetime.
# ---------------------------------------------------------------------------

$mdCommit    = "6a2f38c9282965da156b37e0e7cfd48d0b4e31a1"
$zhCommit    = "dd55eb611e18fcaa5a3c1a301a83ac8ef34e703a"
$deCommit    = "e69f50b3813856463b68562e24796282106814f1"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/"
$zhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/"
$deBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config"

# =====================================================================
# Sheet 1: "Overview"
# =====================================================================
$ws1 = $wb.Worksheets.Item(1)

# Remove pre-existing hyperlinks; they will be fully rebuilt below because
# inserting rows does not relocate the existing hyperlink anchors.
$ws1.Cells.Hyperlinks.Delete()

# Grow the table from 3 data rows to 6 data rows (3 new rows added).
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Hyperlinks.Add($ws1.Range("A2"), ($mdBase + "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md"), "", "", "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md")

$ws1.Range("A3").Value = "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($mdBase + "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md"), "", "", "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md")

$ws1.Range("A4").Value = "807257ea-f579-4b2a-bd51-8b2162483a21.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Hyperlinks.Add($ws1.Range("A4"), ($mdBase + "807257ea-f579-4b2a-bd51-8b2162483a21.md"), "", "", "807257ea-f579-4b2a-bd51-8b2162483a21.md")

$ws1.Range("A5").Value = "a2f32072-c7f3-4270-a7d9-182813d699b8.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Hyperlinks.Add($ws1.Range("A5"), ($mdBase + "a2f32072-c7f3-4270-a7d9-182813d699b8.md"), "", "", "a2f32072-c7f3-4270-a7d9-182813d699b8.md")

$ws1.Range("A6").Value = "a5022a31-c014-47bd-a9e7-232f52e5b19a.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"
$ws1.Hyperlinks.Add($ws1.Range("A6"), ($mdBase + "a5022a31-c014-47bd-a9e7-232f52e5b19a.md"), "", "", "a5022a31-c014-47bd-a9e7-232f52e5b19a.md")

$ws1.Range("A7").Value = ".localization-config"
$ws1.Range("B7").Value = "Not to be localized"
$ws1.Range("C7").Value = "Not to be localized"
$ws1.Hyperlinks.Add($ws1.Range("A7"), $configUrl, "", "", ".localization-config")

# =====================================================================
# Sheet 2: "zh-cn"
# =====================================================================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Hyperlinks.Delete()

$ws2.Rows.Item(2).Insert()
$ws2.Rows.Item(2).Insert()
$ws2.Rows.Item(2).Insert()

# Row 2 : 1a9ac024-...
$ws2.Range("A2").Value = "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "1a9ac024-32a5-4c8a-b8d3-556f1c854616.a1c9686f500289f1cd7423b9f9f3f3e6df72e9ef.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-01 03:21:02"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A2"), ($mdBase + "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md"), "", "", "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), ($zhBase + "1a9ac024-32a5-4c8a-b8d3-556f1c854616.a1c9686f500289f1cd7423b9f9f3f3e6df72e9ef.zh-cn.xlf"), "", "", "1a9ac024-32a5-4c8a-b8d3-556f1c854616.a1c9686f500289f1cd7423b9f9f3f3e6df72e9ef.zh-cn.xlf")

# Row 3 : 24c52df8-...
$ws2.Range("A3").Value = "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.57fdeb8eee6ae139d4b8292e0f6223216460f009.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-01 03:21:02"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($mdBase + "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md"), "", "", "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), ($zhBase + "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.57fdeb8eee6ae139d4b8292e0f6223216460f009.zh-cn.xlf"), "", "", "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.57fdeb8eee6ae139d4b8292e0f6223216460f009.zh-cn.xlf")

# Row 4 : 807257ea-...
$ws2.Range("A4").Value = "807257ea-f579-4b2a-bd51-8b2162483a21.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "807257ea-f579-4b2a-bd51-8b2162483a21.425cf7955d30e283e6058e4feaf365477b718d3f.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-01 03:19:37"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A4"), ($mdBase + "807257ea-f579-4b2a-bd51-8b2162483a21.md"), "", "", "807257ea-f579-4b2a-bd51-8b2162483a21.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), ($zhBase + "807257ea-f579-4b2a-bd51-8b2162483a21.425cf7955d30e283e6058e4feaf365477b718d3f.zh-cn.xlf"), "", "", "807257ea-f579-4b2a-bd51-8b2162483a21.425cf7955d30e283e6058e4feaf365477b718d3f.zh-cn.xlf")

# Row 5 : a2f32072-...
$ws2.Range("A5").Value = "a2f32072-c7f3-4270-a7d9-182813d699b8.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "a2f32072-c7f3-4270-a7d9-182813d699b8.02ba6ff35e88f1734aef3ac7764cb234e6b9fea1.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-03-01 03:21:02"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A5"), ($mdBase + "a2f32072-c7f3-4270-a7d9-182813d699b8.md"), "", "", "a2f32072-c7f3-4270-a7d9-182813d699b8.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), ($zhBase + "a2f32072-c7f3-4270-a7d9-182813d699b8.02ba6ff35e88f1734aef3ac7764cb234e6b9fea1.zh-cn.xlf"), "", "", "a2f32072-c7f3-4270-a7d9-182813d699b8.02ba6ff35e88f1734aef3ac7764cb234e6b9fea1.zh-cn.xlf")

# Row 6 : a5022a31-...
$ws2.Range("A6").Value = "a5022a31-c014-47bd-a9e7-232f52e5b19a.md"
$ws2.Range("B6").Value = "Ready for handoff"
$ws2.Range("C6").Value = "a5022a31-c014-47bd-a9e7-232f52e5b19a.d43064b45d8778adf4cccfe109311a11551617a2.zh-cn.xlf"
$ws2.Range("D6").Value = "2016-03-01 03:19:37"
$ws2.Range("G6").Value = "0001-01-01 00:00:00"
$ws2.Range("H6").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A6"), ($mdBase + "a5022a31-c014-47bd-a9e7-232f52e5b19a.md"), "", "", "a5022a31-c014-47bd-a9e7-232f52e5b19a.md")
$ws2.Hyperlinks.Add($ws2.Range("C6"), ($zhBase + "a5022a31-c014-47bd-a9e7-232f52e5b19a.d43064b45d8778adf4cccfe109311a11551617a2.zh-cn.xlf"), "", "", "a5022a31-c014-47bd-a9e7-232f52e5b19a.d43064b45d8778adf4cccfe109311a11551617a2.zh-cn.xlf")

# Row 7 : .localization-config
$ws2.Range("A7").Value = ".localization-config"
$ws2.Range("B7").Value = "Not to be localized"
$ws2.Range("D7").Value = "0001-01-01 00:00:00"
$ws2.Range("G7").Value = "0001-01-01 00:00:00"
$ws2.Range("H7").Value = "Ignored"
$ws2.Hyperlinks.Add($ws2.Range("A7"), $configUrl, "", "", ".localization-config")

# =====================================================================
# Sheet 3: "de-de"
# =====================================================================
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Hyperlinks.Delete()

$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(2).Insert()

# Row 2 : 1a9ac024-...
$ws3.Range("A2").Value = "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "1a9ac024-32a5-4c8a-b8d3-556f1c854616.a1c9686f500289f1cd7423b9f9f3f3e6df72e9ef.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-01 03:21:14"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A2"), ($mdBase + "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md"), "", "", "1a9ac024-32a5-4c8a-b8d3-556f1c854616.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), ($deBase + "1a9ac024-32a5-4c8a-b8d3-556f1c854616.a1c9686f500289f1cd7423b9f9f3f3e6df72e9ef.de-de.xlf"), "", "", "1a9ac024-32a5-4c8a-b8d3-556f1c854616.a1c9686f500289f1cd7423b9f9f3f3e6df72e9ef.de-de.xlf")

# Row 3 : 24c52df8-...
$ws3.Range("A3").Value = "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.57fdeb8eee6ae139d4b8292e0f6223216460f009.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-01 03:21:14"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($mdBase + "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md"), "", "", "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), ($deBase + "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.57fdeb8eee6ae139d4b8292e0f6223216460f009.de-de.xlf"), "", "", "24c52df8-07ae-4e7d-9cbe-cc774dd37cc6.57fdeb8eee6ae139d4b8292e0f6223216460f009.de-de.xlf")

# Row 4 : 807257ea-...
$ws3.Range("A4").Value = "807257ea-f579-4b2a-bd51-8b2162483a21.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "807257ea-f579-4b2a-bd51-8b2162483a21.425cf7955d30e283e6058e4feaf365477b718d3f.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-01 03:20:31"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A4"), ($mdBase + "807257ea-f579-4b2a-bd51-8b2162483a21.md"), "", "", "807257ea-f579-4b2a-bd51-8b2162483a21.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), ($deBase + "807257ea-f579-4b2a-bd51-8b2162483a21.425cf7955d30e283e6058e4feaf365477b718d3f.de-de.xlf"), "", "", "807257ea-f579-4b2a-bd51-8b2162483a21.425cf7955d30e283e6058e4feaf365477b718d3f.de-de.xlf")

# Row 5 : a2f32072-...
$ws3.Range("A5").Value = "a2f32072-c7f3-4270-a7d9-182813d699b8.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "a2f32072-c7f3-4270-a7d9-182813d699b8.02ba6ff35e88f1734aef3ac7764cb234e6b9fea1.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-01 03:21:14"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A5"), ($mdBase + "a2f32072-c7f3-4270-a7d9-182813d699b8.md"), "", "", "a2f32072-c7f3-4270-a7d9-182813d699b8.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), ($deBase + "a2f32072-c7f3-4270-a7d9-182813d699b8.02ba6ff35e88f1734aef3ac7764cb234e6b9fea1.de-de.xlf"), "", "", "a2f32072-c7f3-4270-a7d9-182813d699b8.02ba6ff35e88f1734aef3ac7764cb234e6b9fea1.de-de.xlf")

# Row 6 : a5022a31-...
$ws3.Range("A6").Value = "a5022a31-c014-47bd-a9e7-232f52e5b19a.md"
$ws3.Range("B6").Value = "Ready for handoff"
$ws3.Range("C6").Value = "a5022a31-c014-47bd-a9e7-232f52e5b19a.d43064b45d8778adf4cccfe109311a11551617a2.de-de.xlf"
$ws3.Range("D6").Value = "2016-03-01 03:20:31"
$ws3.Range("G6").Value = "0001-01-01 00:00:00"
$ws3.Range("H6").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A6"), ($mdBase + "a5022a31-c014-47bd-a9e7-232f52e5b19a.md"), "", "", "a5022a31-c014-47bd-a9e7-232f52e5b19a.md")
$ws3.Hyperlinks.Add($ws3.Range("C6"), ($deBase + "a5022a31-c014-47bd-a9e7-232f52e5b19a.d43064b45d8778adf4cccfe109311a11551617a2.de-de.xlf"), "", "", "a5022a31-c014-47bd-a9e7-232f52e5b19a.d43064b45d8778adf4cccfe109311a11551617a2.de-de.xlf")

# Row 7 : .localization-config
$ws3.Range("A7").Value = ".localization-config"
$ws3.Range("B7").Value = "Not to be localized"
$ws3.Range("D7").Value = "0001-01-01 00:00:00"
$ws3.Range("G7").Value = "0001-01-01 00:00:00"
$ws3.Range("H7").Value = "Ignored"
$ws3.Hyperlinks.Add($ws3.Range("A7"), $configUrl, "", "", ".localization-config")
